# Mark specific vocabulary rows as "Processed" in column C.
# Rows: 11-15, then every 10th row from 37 through 1027.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(11, 12, 13, 14, 15)
for ($r = 37; $r -le 1027; $r += 10) {
    $rows += $r
}

foreach ($r in $rows) {
    $ws.Cells.Item($r, 3).Value = "Processed"
}
